$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 44881
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11250
$ws.Range("O2").Value = 11250
$ws.Range("P2").Value = 11250
$ws.Range("S2").Value = 11250

# Row 3 updates
$ws.Range("D3").Value = 44874
$ws.Range("M3").Value = 200
$ws.Range("P3").Value = 7750
$ws.Range("S3").Value = 7750

# Row 5 updates
$ws.Range("D5").Value = 44923
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 7500
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 7625
$ws.Range("S5").Value = 7625
